# Lesson 18 (3rd ed.) wordlist update: append rows 57-115 (new vocabulary)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(57, 1).Value = 'purpose'
$ws.Cells.Item(57, 2).Value = '目的|もくてき'
$ws.Cells.Item(58, 1).Value = 'eye'
$ws.Cells.Item(58, 2).Value = '目|め'
$ws.Cells.Item(59, 1).Value = 'eye drops'
$ws.Cells.Item(59, 2).Value = '目薬|めぐすり'
$ws.Cells.Item(60, 1).Value = 'the second'
$ws.Cells.Item(60, 2).Value = '二番目|にばんめ'
$ws.Cells.Item(61, 1).Value = 'one''s superiors'
$ws.Cells.Item(61, 2).Value = '目上の人|めうえのひと'
$ws.Cells.Item(62, 1).Value = 'modern'
$ws.Cells.Item(62, 2).Value = '現代的|げんだいてき'
$ws.Cells.Item(63, 1).Value = 'social'
$ws.Cells.Item(63, 2).Value = '社会的|しゃかいてき'
$ws.Cells.Item(64, 1).Value = 'target'
$ws.Cells.Item(64, 2).Value = '的|まと'
$ws.Cells.Item(65, 1).Value = 'Western clothes'
$ws.Cells.Item(65, 2).Value = '洋服|ようふく'
$ws.Cells.Item(66, 1).Value = 'the East'
$ws.Cells.Item(66, 2).Value = '東洋|とうよう'
$ws.Cells.Item(67, 1).Value = 'Western food'
$ws.Cells.Item(67, 2).Value = '洋食|ようしょく'
$ws.Cells.Item(68, 1).Value = 'the Atlantic'
$ws.Cells.Item(68, 2).Value = '大西洋|たいせいよう'
$ws.Cells.Item(69, 1).Value = 'clothes'
$ws.Cells.Item(69, 2).Value = '服|ふく'
$ws.Cells.Item(70, 1).Value = 'Western clothes'
$ws.Cells.Item(70, 2).Value = '洋服|ようふく'
$ws.Cells.Item(71, 1).Value = 'uniform'
$ws.Cells.Item(71, 2).Value = '制服|せいふく'
$ws.Cells.Item(72, 1).Value = 'Japanese clothes'
$ws.Cells.Item(72, 2).Value = '和服|わふく'
$ws.Cells.Item(73, 1).Value = 'cafeteria'
$ws.Cells.Item(73, 2).Value = '食堂|しょくどう'
$ws.Cells.Item(74, 1).Value = 'public hall'
$ws.Cells.Item(74, 2).Value = '公会堂|こうかいどう'
$ws.Cells.Item(75, 1).Value = 'dignified; imposing'
$ws.Cells.Item(75, 2).Value = '堂々とした|どうどうとした'
$ws.Cells.Item(76, 1).Value = 'physical labor'
$ws.Cells.Item(76, 2).Value = '力仕事|ちからしごと'
$ws.Cells.Item(77, 1).Value = 'cooperation'
$ws.Cells.Item(77, 2).Value = '協力|きょうりょく'
$ws.Cells.Item(78, 1).Value = 'endeavor'
$ws.Cells.Item(78, 2).Value = '努力|どりょく'
$ws.Cells.Item(79, 1).Value = 'sumo wrestler'
$ws.Cells.Item(79, 2).Value = '力士|りきし'
$ws.Cells.Item(80, 1).Value = 'class'
$ws.Cells.Item(80, 2).Value = '授業|じゅぎょう'
$ws.Cells.Item(81, 1).Value = 'professor'
$ws.Cells.Item(81, 2).Value = '教授|きょうじゅ'
$ws.Cells.Item(82, 1).Value = 'to be given'
$ws.Cells.Item(82, 2).Value = '授かる|さずかる'
$ws.Cells.Item(83, 1).Value = 'class'
$ws.Cells.Item(83, 2).Value = '授業|じゅぎょう'
$ws.Cells.Item(84, 1).Value = 'occupation'
$ws.Cells.Item(84, 2).Value = '職業|しょくぎょう'
$ws.Cells.Item(85, 1).Value = 'industry'
$ws.Cells.Item(85, 2).Value = '産業|さんぎょう'
$ws.Cells.Item(86, 1).Value = 'service industry'
$ws.Cells.Item(86, 2).Value = 'サービス業|サービスぎょう'
$ws.Cells.Item(87, 1).Value = 'exam'
$ws.Cells.Item(87, 2).Value = '試験|しけん'
$ws.Cells.Item(88, 1).Value = 'game; match'
$ws.Cells.Item(88, 2).Value = '試合|しあい'
$ws.Cells.Item(89, 1).Value = 'entrance exam'
$ws.Cells.Item(89, 2).Value = '入試|にゅうし'
$ws.Cells.Item(90, 1).Value = 'to try'
$ws.Cells.Item(90, 2).Value = '試みる|こころみる'
$ws.Cells.Item(91, 1).Value = 'experiment'
$ws.Cells.Item(91, 2).Value = '実験|じっけん'
$ws.Cells.Item(92, 1).Value = 'experience'
$ws.Cells.Item(92, 2).Value = '経験|けいけん'
$ws.Cells.Item(93, 1).Value = 'taking examination'
$ws.Cells.Item(93, 2).Value = '受験|じゅけん'
$ws.Cells.Item(94, 1).Value = 'to lend'
$ws.Cells.Item(94, 2).Value = '貸す|かす'
$ws.Cells.Item(95, 1).Value = 'lending'
$ws.Cells.Item(95, 2).Value = '貸し出し|かしだし'
$ws.Cells.Item(96, 1).Value = 'rental condo'
$ws.Cells.Item(96, 2).Value = '賃貸マンション|ちんたいマンション'
$ws.Cells.Item(97, 1).Value = 'library'
$ws.Cells.Item(97, 2).Value = '図書館|としょかん'
$ws.Cells.Item(98, 1).Value = 'map'
$ws.Cells.Item(98, 2).Value = '地図|ちず'
$ws.Cells.Item(99, 1).Value = 'figure'
$ws.Cells.Item(99, 2).Value = '図|ず'
$ws.Cells.Item(100, 1).Value = 'signal'
$ws.Cells.Item(100, 2).Value = '合図|あいず'
$ws.Cells.Item(101, 1).Value = 'to attempt'
$ws.Cells.Item(101, 2).Value = '図る|はかる'
$ws.Cells.Item(102, 1).Value = 'Japanese inn'
$ws.Cells.Item(102, 2).Value = '旅館|りょかん'
$ws.Cells.Item(103, 1).Value = 'movie theater'
$ws.Cells.Item(103, 2).Value = '映画館|えいがかん'
$ws.Cells.Item(104, 1).Value = 'embassy'
$ws.Cells.Item(104, 2).Value = '大使館|たいしかん'
$ws.Cells.Item(105, 1).Value = 'to come to an end'
$ws.Cells.Item(105, 2).Value = '終わる|おわる'
$ws.Cells.Item(106, 1).Value = 'end'
$ws.Cells.Item(106, 2).Value = '終わり|おわり'
$ws.Cells.Item(107, 1).Value = 'last stop'
$ws.Cells.Item(107, 2).Value = '終点|しゅうてん'
$ws.Cells.Item(108, 1).Value = 'the last...'
$ws.Cells.Item(108, 2).Value = '最終～|さいしゅう～'
$ws.Cells.Item(109, 1).Value = 'homework'
$ws.Cells.Item(109, 2).Value = '宿題|しゅくだい'
$ws.Cells.Item(110, 1).Value = 'boarding house'
$ws.Cells.Item(110, 2).Value = '下宿|げしゅく'
$ws.Cells.Item(111, 1).Value = 'lodging'
$ws.Cells.Item(111, 2).Value = '宿泊|しゅくはく'
$ws.Cells.Item(112, 1).Value = 'inn'
$ws.Cells.Item(112, 2).Value = '宿|やど'
$ws.Cells.Item(113, 1).Value = 'problem; question'
$ws.Cells.Item(113, 2).Value = '問題|もんだい'
$ws.Cells.Item(114, 1).Value = 'topic of conversation'
$ws.Cells.Item(114, 2).Value = '話題|わだい'
$ws.Cells.Item(115, 1).Value = 'title'
$ws.Cells.Item(115, 2).Value = '題|だい'
